$d = $word.ActiveDocument

# Helper: find the Paragraph object (from $doc.Paragraphs) that contains
# document-position $pos. Using the Paragraphs collection (rather than an
# ad-hoc Document.Range with the same character bounds) is what makes
# Range.Delete() on it actually merge/remove the <w:p> element instead of
# just emptying its text.
function Get-ParagraphAt($doc, $pos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $candidate = $doc.Paragraphs.Item($i)
        if ($candidate.Range.Start -le $pos -and $pos -lt $candidate.Range.End) {
            return $candidate
        }
    }
    return $null
}

# 1) Remove the bold "Externe Partnerschaft:" heading paragraph on the cover
#    page entirely (heading text + its own paragraph mark), which shifts the
#    following placeholder paragraph up into its place.
$headingSearch = $d.Content
$headingFound = $headingSearch.Find.Execute("Externe Partnerschaft:", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($headingFound -eq $true) {
    $headingParagraph = Get-ParagraphAt $d $headingSearch.Start
    if ($headingParagraph -ne $null) {
        [void]$headingParagraph.Range.Delete()
    }
}

# 2) Clear the italic placeholder paragraph's text ("Name der Organisation
#    bzw. des/der Auftraggebenden") but keep the paragraph itself (and its
#    paragraph mark / formatting) intact, leaving an empty italic paragraph.
$placeholderSearch = $d.Content
$placeholderFound = $placeholderSearch.Find.Execute("Name der Organisation bzw. des/der Auftraggebenden", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($placeholderFound -eq $true) {
    $placeholderParagraph = Get-ParagraphAt $d $placeholderSearch.Start
    if ($placeholderParagraph -ne $null) {
        $textOnlyRange = $d.Range($placeholderParagraph.Range.Start, $placeholderParagraph.Range.End - 1)
        $textOnlyRange.Text = ""
    }
}
